$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update teacher1's row (row 2)
$ws.Range("C2").Value = 30
$ws.Range("E2").Value = 1111111111
$ws.Range("F2").Value = "teacher1@gmail.com"

# Row 3 becomes what used to be teacher2's data (row 4), then row 4 is deleted
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "teacher2"
$ws.Range("C3").Value = 35
$ws.Range("D3").Value = "5->Math"
$ws.Range("E3").Value = 2222222222
$ws.Range("F3").Value = "teacher2@gmail.com"
$ws.Range("G3").Value = "2 cairo street"
$ws.Range("H3").Value = $false

# Delete the now-duplicate row 4
$ws.Rows("4:4").Delete()

# Update selection/view
$ws.Range("H8").Select()

# Shift the workbook window's on-screen position (best effort - cosmetic)
try { $wb.Windows.Item(1).Left = 7080 } catch { }
